# Adds two new entity tables ("CARRITOCOMPRAS" and "VENDECARRITO") to the
# schema sheet, below the existing "PROMOCION" table, matching the
# "Se agregó verificación elemento existente en el carrito de compra" commit.
#
# Approach: clone the cell *formatting* (fills/borders/alignment) of the
# already-existing, visually-identical table blocks via Copy/PasteSpecial
# (formats only) so we reuse the workbook's existing style records instead
# of inventing new ones, then fill in the new cell text content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1) CARRITOCOMPRAS table (rows 47-50), formatted like the PROMOCION
#    table above it (rows 41-44).
# ---------------------------------------------------------------------

# Row 47: title row (no explicit per-cell style, just like row 41/A41's
# siblings A21/A36 which are bare title cells) -> leave default format.

# Row 48: header row styling (A:D) cloned from row 42's header row.
$ws.Range("A42:D42").Copy() | Out-Null
$ws.Range("A48:D48").PasteSpecial($xlPasteFormats) | Out-Null

# Row 49: "type" row styling cloned from row 43.
$ws.Range("A43:D43").Copy() | Out-Null
$ws.Range("A49:D49").PasteSpecial($xlPasteFormats) | Out-Null

# Row 50: constraint row styling cloned from row 44 (thick bottom border).
$ws.Range("A44:D44").Copy() | Out-Null
$ws.Range("A50:D50").PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) VENDECARRITO table (rows 53-56), formatted like the VENDE table
#    above it (rows 36-39).
# ---------------------------------------------------------------------

# Row 53: title row -> default format (same pattern as row 36/A36).

# Row 54: header row styling cloned from row 37's left-hand header block.
$ws.Range("A37:C37").Copy() | Out-Null
$ws.Range("A54:C54").PasteSpecial($xlPasteFormats) | Out-Null

# Row 55: "type" row styling cloned from row 38.
$ws.Range("A38:C38").Copy() | Out-Null
$ws.Range("A55:C55").PasteSpecial($xlPasteFormats) | Out-Null

# Row 56: constraint row styling cloned from the B44:D44 block (same
# "s=21" shaded/thick-bottom style used for this table's PK/FK row).
$ws.Range("B44:D44").Copy() | Out-Null
$ws.Range("A56:C56").PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3) Cell text content. Set in the same order the table was authored
#    (titles/headers first, then constraints, then the trailing
#    "cantidadCarrito" column) so new shared-string entries land in a
#    natural sequence.
# ---------------------------------------------------------------------

$ws.Range("A47").Value = "CARRITOCOMPRAS"
$ws.Range("B48").Value = "estado"
$ws.Range("A53").Value = "VENDECARRITO"
$ws.Range("A54").Value = "idCarrito"
$ws.Range("A56").Value = "PK`nFK (Carritocompras.id)"
$ws.Range("B56").Value = "FK (Vende.id)"
$ws.Range("D50").Value = "FK (Vende.id)"
$ws.Range("C54").Value = "cantidadCarrito"

# Remaining cells reuse already-existing shared strings.
$ws.Range("A48").Value = "id"
$ws.Range("C48").Value = "idCliente"
$ws.Range("D48").Value = "idSucursal"

$ws.Range("A49").Value = "Number"
$ws.Range("B49").Value = "Number"
$ws.Range("C49").Value = "Number"
$ws.Range("D49").Value = "Number"

$ws.Range("A50").Value = "PK"
$ws.Range("B50").Value = "NN"

$ws.Range("B54").Value = "idProducto"

$ws.Range("A55").Value = "Number"
$ws.Range("B55").Value = "Number"
$ws.Range("C55").Value = "Number"
